$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 7.955277333333332
$ws.Range("N2").Value = 23.865832
$ws.Range("O2").Value = 0.05015625076675284
$ws.Range("P2").Value = 0.05015625076675283
$ws.Range("Q2").Value = 11.31839734359111
$ws.Range("R2").Value = 101.86557609232
$ws.Range("S2").Value = 0.05015625076675284
$ws.Range("T2").Value = 0.05015625076675283

$ws.Range("M3").Value = 82.48060333333333
$ws.Range("O3").Value = 0.520021823355633
$ws.Range("P3").Value = 0.520021823355633
$ws.Range("S3").Value = 0.520021823355633
$ws.Range("T3").Value = 0.520021823355633

$ws.Range("M4").Value = 66.90297433333333
$ws.Range("N4").Value = 200.708923
$ws.Range("O4").Value = 0.4218083439585467
$ws.Range("P4").Value = 0.4218083439585465
$ws.Range("Q4").Value = 95.18642974266444
$ws.Range("R4").Value = 856.6778676839799
$ws.Range("S4").Value = 0.4218083439585467
$ws.Range("T4").Value = 0.4218083439585465

$ws.Range("M5").Value = 1.271033333333333
$ws.Range("N5").Value = 3.8131
$ws.Range("O5").Value = 0.008013581919067616
$ws.Range("P5").Value = 0.008013581919067614
$ws.Range("Q5").Value = 1.808366911777778
$ws.Range("R5").Value = 16.275302206
$ws.Range("S5").Value = 0.008013581919067616
$ws.Range("T5").Value = 0.008013581919067614
